$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.312.75'
$ws.Range("E2").Value = '  +1.53%  '
$ws.Range("D3").Value = '1.864.16'
$ws.Range("E3").Value = '  +1.20%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.020'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +1.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.020'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.24%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4810'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.26%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07453'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.22%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9373'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.73'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.16%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07879'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.30%  '
$ws.Range("D13").Value = '1.865.86'
$ws.Range("E13").Value = '  +1.54%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.435'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.86%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.551'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.28%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '90.39'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.63%  '
$ws.Range("E17").Value = '  +1.11%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008804'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.06%  '
$ws.Range("E19").Value = '  +1.23%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.82'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.36%  '
$ws.Range("D21").Value = '27.340.91'
$ws.Range("E21").Value = '  +1.51%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.131'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.27%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.71'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.87%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.963'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.19'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.60%  '
$ws.Range("E26").Value = '  +2.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.013'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '115.98'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.49%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.87%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08928'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.23%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.352'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.05%  '
$ws.Range("E32").Value = '  +2.26%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.571'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.16%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7440'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.15%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.675'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.78%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02049'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.48%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.126'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.36%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05297'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.24%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.001'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.17%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5377'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.61%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.129'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.39%  '
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1537'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.66%  '
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.397'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.93%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.67'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.70%  '
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4840'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.96%  '
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.021'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.33%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.688'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.79%  '
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '103.28'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.29%  '
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '66.73'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.55%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06089'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.94%  '
$ws.Range("B51").Value = 'EOS'
$ws.Range("C51").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9011'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.71%  '
